$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value. Every value is written with a leading apostrophe so Excel
# always stores it as literal text (never auto-coerced into a number/date), then
# ClearFormats() strips the quote-prefix/number-format style that the text-coercion
# trick leaves behind, so the cell keeps its original (unstyled) look.
$updates = [ordered]@{
    "D2" = "42.660.47"
    "E2" = "  -1.72%  "
    "D3" = "2.285.44"
    "E3" = "  -3.55%  "
    "D5" = "301.40"
    "E5" = "  -2.78%  "
    "D6" = "97.76"
    "E6" = "  -5.84%  "
    "E7" = "  -1.51%  "
    "E8" = "  +0.02%  "
    "D9" = "0.500"
    "E9" = "  -3.45%  "
    "D10" = "33.65"
    "E10" = "  -5.68%  "
    "B11" = "Dogecoin"
    "C11" = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
    "D11" = "0.0789"
    "E11" = "  -2.05%  "
    "B12" = "OKB"
    "C12" = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
    "D12" = "50.62"
    "E12" = "  -4.90%  "
    "E13" = "  -0.18%  "
    "D14" = "6.66"
    "E14" = "  -3.86%  "
    "D15" = "2.638.66"
    "D16" = "15.30"
    "E16" = "  -1.36%  "
    "D17" = "2.286.04"
    "E17" = "  -3.64%  "
    "D18" = "0.790"
    "E18" = "  -2.24%  "
    "D19" = "42.568.63"
    "E19" = "  -1.90%  "
    "D20" = "0.0₃0896"
    "E20" = "  -1.89%  "
    "D21" = "11.52"
    "E21" = "  -2.92%  "
    "D22" = "6.01"
    "E22" = "  -4.86%  "
    "D23" = "66.72"
    "E23" = "  -2.03%  "
    "D24" = "235.07"
    "E24" = "  -2.14%  "
    "E25" = "  -4.86%  "
    "D26" = "2.49"
    "E26" = "  -4.26%  "
    "E27" = "  +0.19%  "
    "D28" = "24.54"
    "E28" = "  -4.86%  "
    "E29" = "  -0.46%  "
    "D30" = "164.88"
    "E30" = "  +2.07%  "
    "D31" = "33.76"
    "E31" = "  -7.48%  "
    "D32" = "9.11"
    "E32" = "  -3.61%  "
    "E33" = "  +0.00%  "
    "D34" = "4.98"
    "E34" = "  -4.41%  "
    "D35" = "2.40"
    "E35" = "  -4.06%  "
    "D36" = "0.0694"
    "E36" = "  -5.48%  "
    "D37" = "4.38"
    "E37" = "  -6.01%  "
    "D38" = "16.24"
    "E38" = "  -10.40%  "
    "E39" = "  -7.92%  "
    "E40" = "  -7.52%  "
    "E41" = "  -5.16%  "
    "E42" = "  -3.05%  "
    "D43" = "2.41"
    "E43" = "  -7.44%  "
    "D44" = "1.961.91"
    "E44" = "  -3.58%  "
    "E45" = "  -2.39%  "
    "D46" = "17.80"
    "E46" = "  -9.16%  "
    "E47" = "  -7.96%  "
    "E48" = "  -7.93%  "
    "D49" = "53.53"
    "E49" = "  -7.15%  "
    "E50" = "  -3.44%  "
    "B51" = "THORChain"
    "C51" = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
    "D51" = "4.68"
    "E51" = "  -2.82%  "
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    $range.Value = "'" + $updates[$cell]
    $range.ClearFormats()
}
